# Revert "RESTORE" commit: rename AI/ML-specific template text back to the
# generic "Information Technology" template, and restore the blank spacer
# rows that the RESTORE commit had removed.

$wb = $excel.ActiveWorkbook

function Set-BlankRows {
    param($ws, [int[]]$rows)
    foreach ($r in $rows) {
        # Touching OutlineLevel with its own default value (0) is enough to
        # make the engine persist the row as a bare, attribute-less
        # `<row r="N"/>` element without disturbing any existing content.
        $ws.Rows.Item($r).OutlineLevel = 0
    }
}

# ---------------------------------------------------------------------
# Sheet: Instructions & User Guide
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")

$ws1.Range("A1").Value = "Information Technology Comprehensive Budget - User Guide & Instructions"
$ws1.Range("A56").Value = "📋 INFORMATION TECHNOLOGY PROJECT OVERVIEW"
$ws1.Range("B59").Value = "IT Managers, DevOps Engineers, AI Architects, DevOps Engineers..."

Set-BlankRows $ws1 @(2, 10, 20, 28, 37, 45, 54, 55, 60)

# ---------------------------------------------------------------------
# Sheet: Budget Summary
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Budget Summary")

$ws2.Range("A1").Value = "Information Technology - Executive Budget Summary"

Set-BlankRows $ws2 @(2, 6)

# ---------------------------------------------------------------------
# Sheet: Resources
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Resources")

$ws3.Range("A1").Value = "Information Technology - Resources Budget"
$ws3.Range("A4").Value = "IT Managers"
$ws3.Range("A5").Value = "DevOps Engineers"
$ws3.Range("A9").Value = "System Administrators"

Set-BlankRows $ws3 @(2, 11)

# ---------------------------------------------------------------------
# Sheet: Logistics
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Logistics")

$ws4.Range("A1").Value = "Information Technology - Logistics Budget"

Set-BlankRows $ws4 @(2, 9)

# ---------------------------------------------------------------------
# Sheet: Technology
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Technology")

$ws5.Range("A1").Value = "Information Technology - Technology Budget"

Set-BlankRows $ws5 @(2, 10)

# ---------------------------------------------------------------------
# Sheet: Training
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Training")

$ws6.Range("A1").Value = "Information Technology - Training Budget"
$ws6.Range("A4").Value = "IT Certification Programs"

Set-BlankRows $ws6 @(2, 9)

# ---------------------------------------------------------------------
# Sheet: Contingency
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Contingency")

$ws7.Range("A1").Value = "Information Technology - Contingency Budget"

Set-BlankRows $ws7 @(2, 5, 11, 13)

# ---------------------------------------------------------------------
# Sheet: Timeline
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Timeline")

$ws8.Range("A1").Value = "Information Technology - Budget Timeline"

Set-BlankRows $ws8 @(2)
